# Localization status report refresh ("Generate Report for Archive"):
#  - Status moves from "Ready for handoff" to "In Translation" everywhere it
#    appears (Overview zh-cn/de-de summary columns + the per-language Status
#    column on each language sheet).
#  - The Status-related columns are narrower to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update every cell holding the old status text ---------------------
$wsOverview.Range("E2").Value = "In Translation"   # zh-cn status (Overview)
$wsOverview.Range("F2").Value = "In Translation"   # de-de status (Overview)
$wsZhCn.Range("C2").Value     = "In Translation"   # Status column (zh-cn sheet)
$wsDeDe.Range("C2").Value     = "In Translation"   # Status column (de-de sheet)

# --- Narrow the columns that used to size for "Ready for handoff" ------
$newWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # Overview column E
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # Overview column F
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth   # zh-cn column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth   # de-de column C (Status)
